# Adds an "Electrode Locations" column (C) derived from the file name in
# column A (the substring before the first underscore), then sorts all
# data rows (2..last) by that electrode location: first by the leading
# letter(s), then numerically by the trailing number (natural sort, e.g.
# A2, A3, A5, A8, A11, A14, A15, C1, C3, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

# --- Collect existing data (rows 2..lastRow, columns A & B) -----------
$data = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $name = [string]$ws.Cells.Item($r, 1).Value2
    $val  = $ws.Cells.Item($r, 2).Value2

    $loc = $name.Split("_")[0]

    # split the electrode location into its letter prefix and numeric
    # suffix so the rows can be sorted "naturally" (A2 before A11).
    if ($loc -match '^([A-Za-z]+)(\d+)$') {
        $letters = $Matches[1]
        $number  = [int]$Matches[2]
    } else {
        $letters = $loc
        $number  = 0
    }

    $sortKey = $letters + $number.ToString("D4")

    $data += [PSCustomObject]@{
        Name    = $name
        Value   = $val
        Loc     = $loc
        SortKey = $sortKey
    }
}

# --- Sort by (letter prefix, numeric suffix) ---------------------------
$sorted = $data | Sort-Object -Property SortKey

# --- Write header for the new column ------------------------------------
$ws.Cells.Item(1, 3).Value = "Electrode Locations"
$ws.Cells.Item(1, 1).Copy() | Out-Null
$ws.Cells.Item(1, 3).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Write sorted rows back out -----------------------------------------
$r = 2
foreach ($row in $sorted) {
    $ws.Cells.Item($r, 1).Value = $row.Name
    $ws.Cells.Item($r, 2).Value = $row.Value
    $ws.Cells.Item($r, 3).Value = $row.Loc
    $r++
}
